# Edit script: merges a split run on slide 2 and appends two new
# "Title and Content" slides (Go arrays side-note + array limitations).

$p = $ppt.ActivePresentation

# --- 1. Slide 2: merge the 3 runs of the 3rd paragraph in "TextBox 29"
#        into a single run with the combined text. ---
$s2 = $p.Slides.Item(2)
$textBox29 = $s2.Shapes.Item(20)
$tr = $textBox29.TextFrame.TextRange
$para3 = $tr.Paragraphs(3, 1)
# Force a real text change (identical re-assignment is treated as a no-op),
# then set the final merged text so it collapses to one run.
$para3.Text = "__tmp__"
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "return the value at index 4 in our array, which would give us the integer 5 in our case."

# --- 2. Add slide 3: "Side note on Go arrays" (Title and Content layout) ---
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Side note on Go arrays"

$s3body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3body.Text = "Go arrays are values. Array variables represent the whole/entire array. IT IS NOT A POINTER TO THE FIRST ELEMENT IN THE ARRAY (C/C++). "
$null = $s3body.InsertAfter([char]13 + "Passing an array to a function or assigning it to another variable creates a copy of the array elements. ")

# --- 3. Add slide 4: "Limitations on Arrays" (Title and Content layout) ---
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Limitations on Arrays"

$s4body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4body.Text = "Arrays are the most primitive data structure and only have a set fixed length. "
$null = $s4body.InsertAfter([char]13 + "So this leads to a problem, what if we do not know how big of an array we want? We could just create an array that is very large, but that is a waste of computer memory if none of those allocated spots in the array are ever used. We could also create a very small array but what if now we need more size? Copy everything in the current array into another new larger array is costly in time. For example we have an array of 10,000 elements. Now we need to add one more element. We would need to create another array of larger size and copy all 10,000 ")
$null = $s4body.InsertAfter("elements over.")

Write-Output "Slides now: $($p.Slides.Count)"
